$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables (slides 14, 15, 16) with the new
#        built-in table style GUID. Table styles must be changed through
#        Table.ApplyStyle(), not by assigning Table.Style directly.
$newTableStyleId = "{6664FF85-9AD2-4863-BC0D-F2E178E95A39}"
foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}
